$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.610.68'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.933.76'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.34'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.84'
$ws.Range('E6').Value = '  -4.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.547'
$ws.Range('E7').Value = '  -3.93%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').Value = '  -5.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.29'
$ws.Range('E10').Value = '  -5.10%  '
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0842'
$ws.Range('E12').Value = '  -4.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.71'
$ws.Range('E13').Value = '  -4.63%  '
$ws.Range('D14').Value = '3.401.02'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.40'
$ws.Range('E15').Value = '  -5.87%  '
$ws.Range('D16').Value = '2.941.26'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.976'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '51.539.19'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.30'
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.27'
$ws.Range('E20').Value = '  -4.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.13'
$ws.Range('E21').Value = '  -6.76%  '
$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -3.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.81'
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.46'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.67'
$ws.Range('E25').Value = '  -6.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.174'
$ws.Range('E26').Value = '  -7.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.37'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.16'
$ws.Range('E29').Value = '  -6.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.107'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  +2.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.99'
$ws.Range('E32').Value = '  -5.73%  '
$ws.Range('E33').Value = '  -5.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '35.43'
$ws.Range('E34').Value = '  -7.03%  '
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0425'
$ws.Range('E37').Value = '  -4.43%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.83'
$ws.Range('E39').Value = '  +3.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.19'
$ws.Range('E40').Value = '  -6.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.89'
$ws.Range('E41').Value = '  -5.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.114'
$ws.Range('E42').Value = '  -4.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.72'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '120.52'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('E45').Value = '  -1.30%  '
$ws.Range('D46').Value = '2.089.80'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.22'
$ws.Range('E47').Value = '  -8.01%  '
$ws.Range('E48').Value = '  -6.56%  '
$ws.Range('D49').Value = '3.229.42'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('E51').Value = '  -4.93%  '
